$d = $word.ActiveDocument

# --- Step 1: remove the whole "Meta description: ..." paragraph near the top ---
$findRng = $d.Content
$found = $findRng.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Meta description' paragraph"
}
$metaPara = $findRng.Paragraphs(1)
$metaPara.Range.Delete() | Out-Null

# --- Step 2: replace the closing italic "Create a feature image..." paragraph's text
#             with the (former) meta-description sentence, keeping its italic formatting ---
$findRng2 = $d.Content
$found2 = $findRng2.Find.Execute("Create a feature image", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Create a feature image' paragraph"
}
$imgPara = $findRng2.Paragraphs(1)
$imgRng = $imgPara.Range
$imgRng.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark; keep the run's italic formatting
$imgRng.Text = "Experience the underwater world of Atlantis. Play Aquaman slot online for free and enjoy unique mechanics, graphics, and four progressive jackpots."

# --- Step 3: insert a new paragraph, bold, reading the old H1 title text, right
#             before that (now-retargeted) closing paragraph ---
$closingIndex = $d.Paragraphs.Count
$closingPara = $d.Paragraphs($closingIndex)
$closingStart = $closingPara.Range.Start
$insertionPoint = $d.Range($closingStart, $closingStart)

# Insert a placeholder paragraph break (creates a new, empty paragraph right
# before the closing one, at the same index $closingIndex), then overwrite its
# whole (mark-inclusive) range with clean target markup: a leading empty run
# plus a single bold run carrying the text.
$insertionPoint.InsertBefore("X`r")
$newPara = $d.Paragraphs($closingIndex)
$newPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aquaman Slot for Free - Review &amp; Demo 2021</w:t></w:r></w:p>')

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
